$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.436.76'
$ws.Range('E2').Value = '  +0.34%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.868.37'
$ws.Range('E3').Value = '  -0.54%  '

# Row 4
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7063'
$ws.Range('E5').Value = '  -0.64%  '

# Row 6
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.29'
$ws.Range('E6').Value = '  +0.40%  '

# Row 7
$ws.Range('E7').Value = '  -0.04%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07865'
$ws.Range('E8').Value = '  -1.77%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3135'
$ws.Range('E9').Value = '  -0.82%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.51'
$ws.Range('E10').Value = '  -1.93%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07977'
$ws.Range('E11').Value = '  -3.96%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.877.51'
$ws.Range('E12').Value = '  +0.32%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.210'
$ws.Range('E13').Value = '  -0.87%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '93.36'
$ws.Range('E14').Value = '  -1.25%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7000'
$ws.Range('E15').Value = '  -2.03%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.510'
$ws.Range('E16').Value = '  +2.29%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008388'
$ws.Range('E17').Value = '  -1.90%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.452.69'
$ws.Range('E18').Value = '  +0.41%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '252.55'
$ws.Range('E19').Value = '  +3.67%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.133.25'
$ws.Range('E20').Value = '  +0.66%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.12'
$ws.Range('E21').Value = '  -1.34%  '

# Row 22
$ws.Range('E22').Value = '  -0.01%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.620'
$ws.Range('E23').Value = '  -2.39%  '

# Row 24
$ws.Range('E24').Value = '  -0.19%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1555'
$ws.Range('E25').Value = '  -0.37%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.008'

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '161.15'
$ws.Range('E27').Value = '  -1.01%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.73'
$ws.Range('E28').Value = '  +1.04%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.501'
$ws.Range('E29').Value = '  -0.06%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.328'
$ws.Range('E30').Value = '  -2.14%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.267'
$ws.Range('E31').Value = '  -1.31%  '

# Row 32
$ws.Range('E32').Value = '  +1.15%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05311'
$ws.Range('E33').Value = '  -1.36%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.891'
$ws.Range('E34').Value = '  -2.35%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7492'
$ws.Range('E35').Value = '  -2.91%  '

# Row 36
$ws.Range('E36').Value = '  -0.95%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.712'
$ws.Range('E37').Value = '  +1.12%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01883'
$ws.Range('E38').Value = '  -0.20%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.273.40'
$ws.Range('E39').Value = '  +1.04%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.748'
$ws.Range('E40').Value = '  -0.16%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8948'
$ws.Range('E41').Value = '  -1.12%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.092'
$ws.Range('E42').Value = '  -6.10%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '108.83'
$ws.Range('E43').Value = '  -3.71%  '

# Row 44
$ws.Range('E44').Value = '  -3.74%  '

# Row 45
$ws.Range('E45').Value = '  -0.09%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000127'
$ws.Range('E46').Value = '  -3.57%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.033.79'
$ws.Range('E47').Value = '  +0.63%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.579'
$ws.Range('E48').Value = '  +1.32%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.792'

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.5180'
$ws.Range('E50').Value = '  -0.90%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4309'
